$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.080.01'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.644.54'
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').Value = '''217.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = '''0.5202'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.59%  '
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').Value = '''0.06280'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('D10').Value = '''20.42'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').Value = '''0.07759'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '''4.475'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('D13').Value = '1.668.66'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').Value = '1.870.63'
$ws.Range('D15').Value = '''0.5580'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.77%  '
$ws.Range('D16').Value = '0.0₅7992'
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').Value = '''64.76'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '26.075.63'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '''1.005'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Value = '''4.640'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = '''192.40'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').Value = '''10.09'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').Value = '''5.951'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '''146.48'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '''7.164'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').Value = '''0.05614'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.14%  '
$ws.Range('D31').Value = '''1.265'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('D33').Value = '''3.361'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.43%  '
$ws.Range('D34').Value = '''1.597'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('D35').Value = '''2.787'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('D37').Value = '''0.9358'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.49%  '
$ws.Range('D38').Value = '''0.5663'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').Value = '''5.953'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('D40').Value = '''0.01582'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '1.052.37'
$ws.Range('E41').Value = '  -1.45%  '
$ws.Range('D42').Value = '''2.567'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').Value = '''1.005'
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Value = '''0.8425'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('D45').Value = '''102.13'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.70%  '
$ws.Range('D46').Value = '1.782.25'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('D47').Value = '''57.07'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.05325'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₈102'
$ws.Range('E50').Value = '  -3.69%  '
$ws.Range('E51').Value = '  -1.51%  '
